# Applies the "raven.docx" revision described in the commit:
#   1. Append two trailing spaces to the first paragraph's existing text.
#   2. Append a parenthetical note - "(This is a change - Version for main
#      branch)" - in red text, typed/inserted as three separate runs.
#   3. Remove the trailing "ank God almighty, we are free at last." paragraph
#      (the tail half of a word broken across paragraphs at the very end of
#      the document).

$d = $word.ActiveDocument

# --- 1 & 2: first paragraph -------------------------------------------------

$EN_DASH = [char]0x2013

# Append "  " right before the paragraph mark of paragraph 1.
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.SetRange($r1.Start, $r1.End - 1)
$r1.InsertAfter("  ")

# First red run: "(This is a change - Ve"
$p1 = $d.Paragraphs(1)
$r2 = $p1.Range
$r2.SetRange($r2.Start, $r2.End - 1)
$r2.Collapse(0)
$r2.InsertAfter("(This is a change " + $EN_DASH + " Ve")
$r2.Font.Color = 255

# Second red run: "rsion for main branch"
$p1 = $d.Paragraphs(1)
$r3 = $p1.Range
$r3.SetRange($r3.Start, $r3.End - 1)
$r3.Collapse(0)
$r3.InsertAfter("rsion for main branch")
$r3.Font.Color = 255

# Third red run: ")"
$p1 = $d.Paragraphs(1)
$r4 = $p1.Range
$r4.SetRange($r4.Start, $r4.End - 1)
$r4.Collapse(0)
$r4.InsertAfter(")")
$r4.Font.Color = 255

# --- 3: drop the stray "ank God almighty, we are free at last." paragraph --

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)
$lastPara.Range.Delete()
